# Commit message: "decide to remove nonstandard dealer offer"
# Remove the nonstandard dealer discount/premium (F2/H2 = 1 -> 0) and
# adjust the dealer margin threshold table (C9:C16 / D15:D16) to the
# new, standard values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Turn off the nonstandard dealer offer flags ---
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 0

# --- Update the margin/threshold table (rows 9-16) ---
$ws.Range("C9").Value = 330
$ws.Range("C10").Value = 370
$ws.Range("C11").Value = 400
$ws.Range("C12").Value = 410
$ws.Range("C13").Value = 421
$ws.Range("C14").Value = 435
$ws.Range("C15").Value = 460
$ws.Range("D15").Value = 20
$ws.Range("C16").Value = 600
$ws.Range("D16").Value = 30

# --- Restore default view (scroll back to top, select J6) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J6").Select()
